# Apply the FFV schedule ValueSet update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to 2022-01-21T20:46:54+00:00
#  - Publisher value set to "Alvearie Team"
#  - The duplicated "Contact" / "No display for ContactDetail" rows are
#    replaced by a single "Jurisdiction" / "United States of America" row
#    (net: one row removed from the Metadata sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Simple value updates
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"

# Delete the second (duplicate) "Contact" row entirely, shifting rows below up
$ws.Rows.Item(11).Delete()

# Turn the remaining former "Contact" row into the new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
